$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 13.13570827151786
$ws.Range("D2").Value = 4.909754846414764
$ws.Range("E2").Value = 13.38319263054961
$ws.Range("F2").Value = 66.36719626318437
$ws.Range("G2").Value = 91.92550057054049
$ws.Range("H2").Value = 29.17965981741665
$ws.Range("J2").Value = 12.35636330501269
$ws.Range("L2").Value = 9.292654946415571

$ws.Range("C3").Value = 13.13568159928616
$ws.Range("D3").Value = 4.839045006994406
$ws.Range("E3").Value = 13.41966388238335
$ws.Range("F3").Value = 65.86153012665204
$ws.Range("G3").Value = 90.82479694943703
$ws.Range("H3").Value = 29.05688509001493
$ws.Range("J3").Value = 12.38562875396487
$ws.Range("L3").Value = 9.317239672224838

$ws.Range("C4").Value = 13.13886024888697
$ws.Range("D4").Value = 4.79452684284883
$ws.Range("E4").Value = 13.44419776024427
$ws.Range("F4").Value = 65.5711721464409
$ws.Range("G4").Value = 90.17573754448107
$ws.Range("H4").Value = 28.99000105951731
$ws.Range("J4").Value = 12.40586554833595
$ws.Range("L4").Value = 9.333154588727014

$ws.Range("C5").Value = 13.14095568644042
$ws.Range("D5").Value = 4.77611377361212
$ws.Range("E5").Value = 13.45473318117157
$ws.Range("F5").Value = 65.45798041088217
$ws.Range("G5").Value = 89.91822218555492
$ws.Range("H5").Value = 28.96488851197993
$ws.Range("J5").Value = 12.41468085812786
$ws.Range("L5").Value = 9.339847152604168

$ws.Range("C6").Value = 13.14135185321533
$ws.Range("D6").Value = 4.773040062193575
$ws.Range("E6").Value = 13.45651503773448
$ws.Range("F6").Value = 65.43949682095412
$ws.Range("G6").Value = 89.8758905371086
$ws.Range("H6").Value = 28.96084815066814
$ws.Range("J6").Value = 12.41617893244471
$ws.Range("L6").Value = 9.340970984604349

$ws.Range("C7").Value = 13.13888527368342
$ws.Range("D7").Value = 4.794279609809278
$ws.Range("E7").Value = 13.44433766832111
$ws.Range("F7").Value = 65.56962473912996
$ws.Range("G7").Value = 90.17223601880026
$ws.Range("H7").Value = 28.98965369949078
$ws.Range("J7").Value = 12.40598213431083
$ws.Range("L7").Value = 9.333244007138218

$ws.Range("C8").Value = 13.13503406620736
$ws.Range("D8").Value = 4.885606090086728
$ws.Range("E8").Value = 13.39532334705636
$ws.Range("F8").Value = 66.18871196416561
$ws.Range("G8").Value = 91.54057454879062
$ws.Range("H8").Value = 29.13556723852036
$ws.Range("J8").Value = 12.36598242573298
$ws.Range("L8").Value = 9.300962192144571

$ws.Range("C9").Value = 13.15296658096247
$ws.Range("D9").Value = 5.055722423427076
$ws.Range("E9").Value = 13.31621534131701
$ws.Range("F9").Value = 67.55880221309522
$ws.Range("G9").Value = 94.4246634736959
$ws.Range("H9").Value = 29.48876504272335
$ws.Range("J9").Value = 12.30560763148374
$ws.Range("L9").Value = 9.244118589551075

$ws.Range("C10").Value = 13.1818401774222
$ws.Range("D10").Value = 5.174948579168842
$ws.Range("E10").Value = 13.26850148064438
$ws.Range("F10").Value = 68.65511171832938
$ws.Range("G10").Value = 96.64949959344
$ws.Range("H10").Value = 29.78845194623609
$ws.Range("J10").Value = 12.2723620644891
$ws.Range("L10").Value = 9.20623321114115

$ws.Range("C11").Value = 13.19841109087208
$ws.Range("D11").Value = 5.227881757013651
$ws.Range("E11").Value = 13.24906354478732
$ws.Range("F11").Value = 69.17200953778018
$ws.Range("G11").Value = 97.68094682909813
$ws.Range("H11").Value = 29.93331183983026
$ws.Range("J11").Value = 12.25967186736832
$ws.Range("L11").Value = 9.189826973245269

$ws.Range("C12").Value = 13.20518160696991
$ws.Range("D12").Value = 5.247734323332162
$ws.Range("E12").Value = 13.24202969478385
$ws.Range("F12").Value = 69.37023878154565
$ws.Range("G12").Value = 98.07399808551641
$ws.Range("H12").Value = 29.98937197065198
$ws.Range("J12").Value = 12.25521805012022
$ws.Range("L12").Value = 9.183732382461793

$ws.Range("C13").Value = 13.2037014018067
$ws.Range("D13").Value = 5.243467331644001
$ws.Range("E13").Value = 13.24353000889991
$ws.Range("F13").Value = 69.32743778010814
$ws.Range("G13").Value = 97.98924292881209
$ws.Range("H13").Value = 29.97724520575045
$ws.Range("J13").Value = 12.2561615900106
$ws.Range("L13").Value = 9.185039724830938

$ws.Range("C14").Value = 13.19895816903233
$ws.Range("D14").Value = 5.229518914718803
$ws.Range("E14").Value = 13.2484783098133
$ws.Range("F14").Value = 69.18826877130307
$ws.Range("G14").Value = 97.71323582550517
$ws.Range("H14").Value = 29.93789992741315
$ws.Range("J14").Value = 12.25929839075901
$ws.Range("L14").Value = 9.189323206371029

$ws.Range("C15").Value = 13.19611736010354
$ws.Range("D15").Value = 5.220949959231278
$ws.Range("E15").Value = 13.25155188081934
$ws.Range("F15").Value = 69.10334443098397
$ws.Range("G15").Value = 97.54448499807039
$ws.Range("H15").Value = 29.91395602455483
$ws.Range("J15").Value = 12.26126562388254
$ws.Range("L15").Value = 9.191962313089464

$ws.Range("C16").Value = 13.18082649889564
$ws.Range("D16").Value = 5.171462639368317
$ws.Range("E16").Value = 13.26981749091238
$ws.Range("F16").Value = 68.6216866681157
$ws.Range("G16").Value = 96.58245435307131
$ws.Range("H16").Value = 29.77915485500474
$ws.Range("J16").Value = 12.27324051585295
$ws.Range("L16").Value = 9.207321976127785

$ws.Range("C17").Value = 13.1723272302077
$ws.Range("D17").Value = 5.140766536580003
$ws.Range("E17").Value = 13.28160415744337
$ws.Range("F17").Value = 68.3307758736418
$ws.Range("G17").Value = 95.99701092761076
$ws.Range("H17").Value = 29.69862881748112
$ws.Range("J17").Value = 12.28121123155868
$ws.Range("L17").Value = 9.216956015557567

$ws.Range("C18").Value = 13.16776204962487
$ws.Range("D18").Value = 5.122988527992419
$ws.Range("E18").Value = 13.28859689783854
$ws.Range("F18").Value = 68.16517070828479
$ws.Range("G18").Value = 95.66212064698635
$ws.Range("H18").Value = 29.6531167959368
$ws.Range("J18").Value = 12.28602471241932
$ws.Range("L18").Value = 9.222575254907655

$ws.Range("C19").Value = 13.16627184971428
$ws.Range("D19").Value = 5.11694832029606
$ws.Range("E19").Value = 13.29100113954185
$ws.Range("F19").Value = 68.10939855761877
$ws.Range("G19").Value = 95.54905851210937
$ws.Range("H19").Value = 29.63784599671751
$ws.Range("J19").Value = 12.28769374024933
$ws.Range("L19").Value = 9.224491257960455

$ws.Range("C20").Value = 13.17319851278442
$ws.Range("D20").Value = 5.144046893982062
$ws.Range("E20").Value = 13.28032735859436
$ws.Range("F20").Value = 68.36156678764937
$ws.Range("G20").Value = 96.05914412413227
$ws.Range("H20").Value = 29.70711784796438
$ws.Range("J20").Value = 12.28033902842985
$ws.Range("L20").Value = 9.215922390486037

$ws.Range("C21").Value = 13.20033791579092
$ws.Range("D21").Value = 5.233621151473041
$ws.Range("E21").Value = 13.24701599585174
$ws.Range("F21").Value = 69.22907947623965
$ws.Range("G21").Value = 97.7942414559444
$ws.Range("H21").Value = 29.94942407620647
$ws.Range("J21").Value = 12.25836747766041
$ws.Range("L21").Value = 9.188061846697062

$ws.Range("C22").Value = 13.22096340698393
$ws.Range("D22").Value = 5.291041576425811
$ws.Range("E22").Value = 13.22715077172116
$ws.Range("F22").Value = 69.81050981515375
$ws.Range("G22").Value = 98.94244608913321
$ws.Range("H22").Value = 30.11479618349566
$ws.Range("J22").Value = 12.24605850294754
$ws.Range("L22").Value = 9.17054126698886

$ws.Range("C23").Value = 13.20969057098997
$ws.Range("D23").Value = 5.260499294461621
$ws.Range("E23").Value = 13.23757856487296
$ws.Range("F23").Value = 69.49890823918574
$ws.Range("G23").Value = 98.32843190591927
$ws.Range("H23").Value = 30.02590027946237
$ws.Range("J23").Value = 12.25243981281
$ws.Range("L23").Value = 9.17982970983503

$ws.Range("C24").Value = 13.17280360569022
$ws.Range("D24").Value = 5.142564250242425
$ws.Range("E24").Value = 13.28090392528961
$ws.Range("F24").Value = 68.34764108755085
$ws.Range("G24").Value = 96.0310484210961
$ws.Range("H24").Value = 29.70327751531577
$ws.Range("J24").Value = 12.28073263205317
$ws.Range("L24").Value = 9.216389441562507

$ws.Range("C25").Value = 13.14536869356386
$ws.Range("D25").Value = 5.010692816027417
$ws.Range("E25").Value = 13.33579140500631
$ws.Range("F25").Value = 67.17197704623464
$ws.Range("G25").Value = 93.6246110138855
$ws.Range("H25").Value = 29.3860965097931
$ws.Range("J25").Value = 12.31999639875616
$ws.Range("L25").Value = 9.258810910984709
